# Add two new rows (7 and 8) of statistics data to the "Statistic" sheet,
# mirroring the structure of the existing rows (in particular row 5, which
# holds a "run terminated early" style record: a handful of populated
# numeric columns, a run of columns that are present but blank, and a
# trailing zero in column AC).

try {
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns (1-based) that stay blank (typed as empty text) for both new rows,
# mirroring row 5's K..AB run of empty-string cells.
$blankCols = 11..28   # K .. AB

$rows = 7, 8
$aVals = 5, 6

for ($i = 0; $i -lt 2; $i++) {
    $r = $rows[$i]
    $aVal = $aVals[$i]

    $ws.Cells.Item($r, 1).Value  = $aVal        # A
    $ws.Cells.Item($r, 2).Value  = 194700       # B
    $ws.Cells.Item($r, 3).Value  = 251          # C
    $ws.Cells.Item($r, 4).Value  = 2            # D
    $ws.Cells.Item($r, 5).Value  = 2            # E
    # F intentionally left empty (matches source row's gap)
    $ws.Cells.Item($r, 7).Value  = 0.00796812749003984   # G
    $ws.Cells.Item($r, 8).Value  = 1                     # H
    # I intentionally left empty
    $ws.Cells.Item($r, 10).Value = 0.9920318725099602    # J

    foreach ($col in $blankCols) {
        # A formula literal of "" is the only way, through the Excel object
        # model, to persist a cell that is typed as a (blank) string rather
        # than leaving the cell absent entirely.
        $ws.Cells.Item($r, $col).Formula = "=""" + """"
    }

    $ws.Cells.Item($r, 29).Value = 0   # AC
}

Write-Output "Added rows 7 and 8 to Statistic sheet"
} catch {
Write-Output "ERROR:"
Write-Output $_.Exception.Message
}
